# Rename the two exercise title slides so both reflect they are part of
# "Exercise 1" (part 1 and part 2), splitting each title into two runs:
# "Exercise " + "1 part N", matching the original author's edit.

$p = $ppt.ActivePresentation

# Slide 10: title "Exercise 1" -> "Exercise " / "1 part 1"
$s10 = $p.Slides.Item(10)
$tr10 = $s10.Shapes.Item(1).TextFrame.TextRange
$tr10.Text = "Exercise "
[void]$tr10.InsertAfter("1 part 1")

# Slide 11: title "Exercise 2" -> "Exercise " / "1 part 2"
$s11 = $p.Slides.Item(11)
$tr11 = $s11.Shapes.Item(1).TextFrame.TextRange
$tr11.Text = "Exercise "
[void]$tr11.InsertAfter("1 part 2")
